$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for A1 ("Tiefe " -> "Tiefe in mysec")
$ws.Range("A1").Value = "Tiefe in mysec"

# Replace the literal/leftover values in column B (rows 2-17) with a formula
# that derives the time-of-flight value from column A: (Tiefe-12.28)*6/4
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Formula = "=(A$r-12.28)*6/4"
    # keep the default (unstyled) formatting, matching the neighbouring column C
    $cell.Style = $ws.Cells.Item($r, 3).Style
}

# Move/restore the active cell selection to B17
[void]$ws.Range("B17").Select()
